$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Saved on" timestamp in the study description (A2)
$ws.Range("A2").Value = "This is an interesting study. Saved on : 2022/09/26 17:01:55"

# Highlight statistically significant p-values (p < 0.05) in the
# OR (univariate) / OR (model 1) / OR (model 2) columns (E, F, G) with a
# light-blue fill and left/center alignment, as the function added in this
# commit does when building the save-summary table.
$highlightColor = 15128749   # RGB(173, 216, 230) = ADD8E6, light blue
$highlightPatternColor = 16443110  # RGB(230, 230, 250) = E6E6FA, lavender

$significantCells = @("E5", "F5", "G5", "E6", "F6", "G6", "E8", "F8", "E9", "F9", "E10", "F10", "G10")

foreach ($cellRef in $significantCells) {
    $cell = $ws.Range($cellRef)
    $cell.HorizontalAlignment = -4131   # xlLeft
    $cell.VerticalAlignment = -4108     # xlCenter
    $cell.Interior.Color = $highlightColor
    $cell.Interior.PatternColor = $highlightPatternColor
}
